# Ansprache wird jetzt unterstützt
# 1) Body text: merge the split "(Jer 29,11-14a)" runs (which were broken up
#    around a spell-check proofErr marker on "Jer") back into a single run.
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "(Jer 29,11-14a)", $false, $false, $false, $false, $false,
    $true, 1, $false, "(Jer 29,11-14a)", 2) | Out-Null

# 2) Default header keyword list: swap in the new picture-description tags.
$header = $d.Sections(1).Headers(1)
$header.Range.Find.Execute(
    "Weg, Weg, Freude, Psalm23, Säulen, Bild", $false, $false, $false, $false, $false,
    $true, 1, $false, "Fluss, Weg, Regenbogen, Bild, Hand, Säulen", 2) | Out-Null
